$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M")

$data = @{
    2 = @(1.019999999999999, 1.057028082059764, 1.054582980607437, 1.06229614213626, 1.071272194292041, 1.041541081105195, 1.062026641344852, 1.057325791949646, 1.065017874336167, 1.073969809336351)
    3 = @(1.02, 1.058887323096292, 1.056009820986789, 1.064012694455673, 1.073179524846739, 1.042037246872855, 1.063533812675358, 1.058564337102249, 1.06654694832461, 1.075690969176053)
    4 = @(1.02, 1.060086506140558, 1.056929563082408, 1.065119628164645, 1.074410140828673, 1.04235504896253, 1.064504932987289, 1.059361726369347, 1.067532116180818, 1.076800699586078)
    5 = @(1.02, 1.060589736302968, 1.057315394545963, 1.065584094141732, 1.074926661026742, 1.042487881101728, 1.064912222154091, 1.059695997014648, 1.067945281716253, 1.07726629740958)
    6 = @(1.02, 1.060674178291913, 1.057380129199358, 1.065662028427668, 1.075013338906869, 1.042510139121103, 1.064980551405993, 1.059752067102849, 1.068014595942984, 1.077344419173303)
    7 = @(1.02, 1.060093233860751, 1.056934721815671, 1.065125837849638, 1.074417045833764, 1.042356826896374, 1.0645103789932, 1.059366196631577, 1.067537640820941, 1.076806924566867)
    8 = @(1.02, 1.05765723482142, 1.055065925320917, 1.062877055041987, 1.071917533755745, 1.041709439456033, 1.062536860137194, 1.057745207179451, 1.065535521520612, 1.074552317809331)
    9 = @(1.02, 1.053334141898693, 1.051745266633999, 1.058884547835645, 1.067484940558738, 1.040543482043444, 1.059026961010733, 1.054857322554998, 1.061974236987057, 1.070548118992964)
    10 = @(1.02, 1.050430258943209, 1.049512023025341, 1.056201623539915, 1.064509682149195, 1.03974884742451, 1.056664251418192, 1.052910029353236, 1.059576563544745, 1.067856418751998)
    11 = @(1.02, 1.049167382954213, 1.048540180199645, 1.055034583741953, 1.063216285067726, 1.039400561045358, 1.055635532853797, 1.052061409389597, 1.058532528170123, 1.066685341407286)
    12 = @(1.02, 1.04869744707706, 1.048178450458765, 1.054600271743848, 1.062735070133364, 1.039270552420888, 1.055252551610408, 1.051745361641217, 1.05814383015294, 1.066249493596748)
    13 = @(1.02, 1.04879828867228, 1.048256076674938, 1.054693470683709, 1.062838328507407, 1.039298468781499, 1.055334742060853, 1.051813192901789, 1.058227248124881, 1.066343023613213)
    14 = @(1.02, 1.049128555341979, 1.048510294788876, 1.054998700267171, 1.063176523945253, 1.039389827584781, 1.055603893359488, 1.052035301865346, 1.058500416689692, 1.066649331720319)
    15 = @(1.02, 1.049331930484706, 1.048666827827674, 1.055186652749111, 1.063384791927597, 1.039446031772953, 1.055769610538446, 1.05217203970484, 1.058668605380344, 1.066837943922044)
    16 = @(1.02, 1.050513954385574, 1.049576417676735, 1.056278962188932, 1.064595411286394, 1.039771872792072, 1.056732403299951, 1.052966233647403, 1.059645728150368, 1.067934020321292)
    17 = @(1.02, 1.051253923349349, 1.050145672743244, 1.05696269982927, 1.065353420486848, 1.039975133091841, 1.057334810748346, 1.053462945173145, 1.060257076419968, 1.068620057004895)
    18 = @(1.02, 1.051685008116541, 1.050477244438253, 1.057361001077914, 1.065795064960651, 1.040093286267254, 1.057685641107993, 1.053752145917171, 1.060613104763189, 1.069019676358537)
    19 = @(1.02, 1.051831908307284, 1.050590223455658, 1.057496725462484, 1.065945572038393, 1.040133504968791, 1.057805173588995, 1.053850667746662, 1.060734406755159, 1.069155846312812)
    20 = @(1.02, 1.051174586331477, 1.050084645303678, 1.056889394306839, 1.065272144080575, 1.039953367145538, 1.057270234517173, 1.053409706912371, 1.060191542692244, 1.068546507187212)
    21 = @(1.02, 1.049031323636169, 1.048435454574414, 1.054908840623337, 1.063076955818271, 1.039362942413179, 1.055524659175707, 1.051969919429258, 1.058420000284767, 1.066559155477641)
    22 = @(1.02, 1.047678853340303, 1.047394229069326, 1.0536588237096, 1.061692172479676, 1.038988014677525, 1.054422106407931, 1.05105984273556, 1.057300961899153, 1.065304655376599)
    23 = @(1.02, 1.048396297804225, 1.047946617337867, 1.05432194061251, 1.062426714979057, 1.039187124754817, 1.055007075188574, 1.051542754664305, 1.057894685366543, 1.065970168849648)
    24 = @(1.02, 1.051210436941367, 1.050112222413634, 1.056922519496479, 1.065308870901851, 1.039963203491385, 1.057299415415909, 1.053433764605471, 1.060221156295085, 1.068579742834062)
    25 = @(1.02, 1.054455516298978, 1.052607100643386, 1.059920365656841, 1.06863433698814, 1.040847934668216, 1.059938290402981, 1.055607730173483, 1.062898972191632, 1.071587130839684)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

Write-Host "Done updating vm_pu values"